$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.842.58'
$ws.Range("E2").Value = '  -1.46%  '
# Row 3
$ws.Range("D3").Value = '2.681.73'
$ws.Range("E3").Value = '  -1.96%  '
# Row 4
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.83'
$ws.Range("E5").Value = '  -1.82%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.03'
$ws.Range("E6").Value = '  -0.78%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.00%  '
# Row 8
$ws.Range("E8").Value = '  -0.93%  '
# Row 9
$ws.Range("E9").Value = '  -2.77%  '
# Row 10
$ws.Range("E10").Value = '  -2.04%  '
# Row 11
$ws.Range("E11").Value = '  -2.80%  '
# Row 12
$ws.Range("E12").Value = '  -3.72%  '
# Row 13
$ws.Range("D13").Value = '3.155.09'
$ws.Range("E13").Value = '  -1.97%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.53'
$ws.Range("E14").Value = '  -1.46%  '
# Row 15
$ws.Range("D15").Value = '62.760.43'
# Row 16
$ws.Range("E16").Value = '  -1.67%  '
# Row 17
$ws.Range("D17").Value = '2.681.61'
$ws.Range("E17").Value = '  -2.13%  '
# Row 18
$ws.Range("E18").Value = '  -3.71%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.61'
$ws.Range("E19").Value = '  -2.65%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '345.06'
$ws.Range("E20").Value = '  -2.26%  '
# Row 21
$ws.Range("E21").Value = '  -4.29%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.07%  '
# Row 23
$ws.Range("E23").Value = '  -2.81%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.13'
$ws.Range("E24").Value = '  -1.66%  '
# Row 25
$ws.Range("E25").Value = '  -0.11%  '
# Row 26
$ws.Range("E26").Value = '  -0.07%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.18'
$ws.Range("E27").Value = '  -1.88%  '
# Row 28
$ws.Range("E28").Value = '  +8.90%  '
# Row 29
$ws.Range("D29").Value = '0.0₃0857'
$ws.Range("E29").Value = '  -5.17%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  +0.46%  '
# Row 31
$ws.Range("E31").Value = '  -1.02%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.19'
$ws.Range("E32").Value = '  +0.77%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.92'
$ws.Range("E33").Value = '  +0.96%  '
# Row 34
$ws.Range("E34").Value = '  +0.56%  '
# Row 35
$ws.Range("E35").Value = '  +0.02%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.47'
$ws.Range("E36").Value = '  -2.73%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("E37").Value = '  +0.29%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '349.49'
$ws.Range("E38").Value = '  +1.63%  '
# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.25'
$ws.Range("E39").Value = '  -0.34%  '
# Row 40
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.947'
$ws.Range("E40").Value = '  -3.14%  '
# Row 41
$ws.Range("E41").Value = '  -2.12%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.33'
$ws.Range("E42").Value = '  -0.21%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.91'
$ws.Range("E43").Value = '  -3.83%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.19'
$ws.Range("E44").Value = '  -3.41%  '
# Row 45
$ws.Range("E45").Value = '  -0.83%  '
# Row 46
$ws.Range("E46").Value = '  -3.56%  '
# Row 47
$ws.Range("E47").Value = '  +0.02%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.00'
$ws.Range("E48").Value = '  -0.67%  '
# Row 49
$ws.Range("E49").Value = '  -2.93%  '
# Row 50
$ws.Range("E50").Value = '  -2.96%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '128.78'
$ws.Range("E51").Value = '  -4.11%  '
